# Generate Report for Handback
# Updates the "zh-cn" and "de-de" localization-status sheets: row 8 (the
# 61762564-d2f3-403a-8f00-616fca4064c8 entry) now has a stale handback
# warning recorded in the "Latest Target File" / "Latest Handback DateTime"
# / "Error Detail" columns, plus a new hyperlink on the target-file cell.
# The Error Detail column is also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff9b89e4fef0346c25fc767f39e4804bfdb9264b/e2e/61762564-d2f3-403a-8f00-616fca4064c8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db5571067f8128555fb9a5254ad2fd7e7796c44f/e2e/61762564-d2f3-403a-8f00-616fca4064c8.md."
$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db5571067f8128555fb9a5254ad2fd7e7796c44f/e2e/61762564-d2f3-403a-8f00-616fca4064c8.md"

$sheetDates = @{ "zh-cn" = "2016-08-22 12:45:15"; "de-de" = "2016-08-22 12:45:24" }

foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # Error Detail column (P) is now wide enough to read the message.
    $ws.Columns.Item(16).ColumnWidth = 40

    $sourceName = $ws.Range("A8").Value2
    $latestTargetFile = $ws.Range("G8").Value2

    # Latest Target File (I8) mirrors the source file name and becomes a
    # hyperlink, same as the existing I2 row for the first entry.
    $ws.Range("I8").Value = $sourceName
    $ws.Hyperlinks.Add($ws.Range("I8"), $latestHandbackUrl, "", "", $sourceName)

    # Latest Handback File (J8) now records the generated xlf file name.
    $ws.Range("J8").Value = $latestTargetFile

    # Latest Handback DateTime (K8).
    $ws.Range("K8").Value = $sheetDates[$name]

    # Error Detail (P8) records the stale-handback-version warning.
    $ws.Range("P8").Value = $errorDetail
}
